$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dll1"
$ws.Range("C2").Value = "Notch1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.190862666666667
$ws.Range("H2").Value = 15.572588
$ws.Range("I2").Value = 0.8740249884703439
$ws.Range("J2").Value = 0.874024988470344
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 45.95651366666667
$ws.Range("N2").Value = 137.869541
$ws.Range("O2").Value = 0.6189188856627118
$ws.Range("P2").Value = 0.6189188856627118
$ws.Range("Q2").Value = 238.5539510824565
$ws.Range("R2").Value = 2146.985559742108
$ws.Range("S2").Value = 0.5409505719054298
$ws.Range("T2").Value = 0.5409505719054298

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dll1"
$ws.Range("C3").Value = "Notch1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.190862666666667
$ws.Range("H3").Value = 15.572588
$ws.Range("I3").Value = 0.8740249884703439
$ws.Range("J3").Value = 0.874024988470344
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.09225115688993263
$ws.Range("P3").Value = 0.09225115688993261
$ws.Range("Q3").Value = 35.55696631305245
$ws.Range("R3").Value = 320.012696817472
$ws.Range("S3").Value = 0.08062981633709926
$ws.Range("T3").Value = 0.08062981633709926

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dll1"
$ws.Range("C4").Value = "Notch1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.190862666666667
$ws.Range("H4").Value = 15.572588
$ws.Range("I4").Value = 0.8740249884703439
$ws.Range("J4").Value = 0.874024988470344
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.446458
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2888299574473556
$ws.Range("P4").Value = 0.2888299574473556
$ws.Range("Q4").Value = 111.3256181644346
$ws.Range("R4").Value = 1001.930563479912
$ws.Range("S4").Value = 0.252444600227815
$ws.Range("T4").Value = 0.2524446002278149

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dll1"
$ws.Range("C5").Value = "Notch1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.100996
$ws.Range("H5").Value = 0.302988
$ws.Range("I5").Value = 0.0170054639091879
$ws.Range("J5").Value = 0.0170054639091879
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 45.95651366666667
$ws.Range("N5").Value = 137.869541
$ws.Range("O5").Value = 0.6189188856627118
$ws.Range("P5").Value = 0.6189188856627118
$ws.Range("Q5").Value = 4.641424054278666
$ws.Range("R5").Value = 41.772816488508
$ws.Range("S5").Value = 0.01052500277285204
$ws.Range("T5").Value = 0.01052500277285204

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dll1"
$ws.Range("C6").Value = "Notch1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.100996
$ws.Range("H6").Value = 0.302988
$ws.Range("I6").Value = 0.0170054639091879
$ws.Range("J6").Value = 0.0170054639091879
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("N6").Value = 20.549744
$ws.Range("O6").Value = 0.09225115688993263
$ws.Range("P6").Value = 0.09225115688993261
$ws.Range("Q6").Value = 0.6918139816746666
$ws.Range("R6").Value = 6.226325835072
$ws.Range("S6").Value = 0.00156877371907258
$ws.Range("T6").Value = 0.00156877371907258

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dll1"
$ws.Range("C7").Value = "Notch1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.100996
$ws.Range("H7").Value = 0.302988
$ws.Range("I7").Value = 0.0170054639091879
$ws.Range("J7").Value = 0.0170054639091879
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.446458
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2888299574473556
$ws.Range("P7").Value = 0.2888299574473556
$ws.Range("Q7").Value = 2.166006472167999
$ws.Range("R7").Value = 19.494058249512
$ws.Range("S7").Value = 0.004911687417263282
$ws.Range("T7").Value = 0.004911687417263281

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Dll1"
$ws.Range("C8").Value = "Notch1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6471736666666666
$ws.Range("H8").Value = 1.941521
$ws.Range("I8").Value = 0.1089695476204681
$ws.Range("J8").Value = 0.1089695476204681
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 45.95651366666667
$ws.Range("N8").Value = 137.869541
$ws.Range("O8").Value = 0.6189188856627118
$ws.Range("P8").Value = 0.6189188856627118
$ws.Range("Q8").Value = 29.74184545687344
$ws.Range("R8").Value = 267.676609111861
$ws.Range("S8").Value = 0.06744331098442995
$ws.Range("T8").Value = 0.06744331098442995

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Dll1"
$ws.Range("C9").Value = "Notch1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6471736666666666
$ws.Range("H9").Value = 1.941521
$ws.Range("I9").Value = 0.1089695476204681
$ws.Range("J9").Value = 0.1089695476204681
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("N9").Value = 20.549744
$ws.Range("O9").Value = 0.09225115688993263
$ws.Range("P9").Value = 0.09225115688993261
$ws.Range("Q9").Value = 4.433084391180445
$ws.Range("R9").Value = 39.897759520624
$ws.Range("S9").Value = 0.01005256683376079
$ws.Range("T9").Value = 0.01005256683376079

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Dll1"
$ws.Range("C10").Value = "Notch1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6471736666666666
$ws.Range("H10").Value = 1.941521
$ws.Range("I10").Value = 0.1089695476204681
$ws.Range("J10").Value = 0.1089695476204681
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.446458
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2888299574473556
$ws.Range("P10").Value = 0.2888299574473556
$ws.Range("Q10").Value = 13.87958286087266
$ws.Range("R10").Value = 124.916245747854
$ws.Range("S10").Value = 0.03147366980227741
$ws.Range("T10").Value = 0.0314736698022774
